$d = $word.ActiveDocument
$rsq = [char]0x2019

# ---------------------------------------------------------------------------
# 1) Expand the "Info:" copy for the DashUI project (matthewia/projects/4).
# ---------------------------------------------------------------------------
$oldInfo = "For Team Sunergy, Appalachian State University" + $rsq + "s Solar Vehicle Team, I designed a digital dashboard interface for their vehicle."
$newInfo = "For Team Sunergy, Appalachian State University" + $rsq + "s Solar Vehicle Team, I designed a digital dashboard interface for their 2018 vehicle, ROSE. We wanted to have a dashboard that evoked the feeling of being in a modern vehicle, but with additional tools fit for a solar-powered electric vehicle (EV). The default screen needed to be easy for the driver to use, but also provide alternate views for testing purposes. "

$d.Content.Find.Execute($oldInfo, $true, $false, $false, $false, $false, $true, 1, $false, $newInfo, 2) | Out-Null

# ---------------------------------------------------------------------------
# 2) Insert a brand new "Copy:" section for DashUI, right after the "Media:"
#    paragraph and before the blank separator that precedes the
#    "matthewia/projects/5" heading.
# ---------------------------------------------------------------------------

# Locate the DashUI "Media:" paragraph by scanning the paragraph collection
# directly (more reliable than combining Find with the Paragraphs collection).
$mediaPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $cand = $d.Paragraphs.Item($i)
    if ($cand.Range.Text -like "Media:*Short clip/video of navigation*") {
        $mediaPara = $cand
        break
    }
}

$anchorPara = $mediaPara.Next()     # the blank spacer paragraph right after "Media:"

$copyText = "Copy:`t1) Not being a developer on this project allowed me to focus my effort in the design. I began by researching existing digital dashboard designs, particularly in EVs such as Teslas. From there, I started sketching wireframes ranging from relatively traditional twin-dial layouts to those befitting a sci-fi film.`r" + `
"`r" + `
"2) I had numerous meetings with the telemetry lead and various directors from the sub-teams, in order to determine what data values are most necessary to display for the driver. In solar vehicle racing introduces a new set of measurements to be monitoring compared to a standard EV, let alone a traditional gasoline car. However, I wanted the design to be approachable to the average user and avoid overwhelming them with too many dials or moving parts. To accomplish that I focused on designing the layout while thinking about the hierarchy of priority of each data field. `r" + `
"`r" + `
"3) A significant element of the design is the large colored arch across the bottom half of the Standard view. This, paired with the value displayed under the arch, was my solution to including the Net Power as both a raw value but also as a quick visual indicator. Net Power is important to solar racing, as it correlates to current efficiency based on your incoming power from the solar panels, and your outgoing power from the motors.`r" + `
"`r" + `
"4) Another aspect of this project required me to effectively communicate my design to the developer working on building the interface. I created a design handoff for the developer to use a reference, as well as had meetings about implementing the design in code.`r" + `
"`r" + `
"`r"

$insertRange = $anchorPara.Range
$insertRange.Collapse(0)
$insertRange.InsertBefore($copyText)

# ---------------------------------------------------------------------------
# 3) Re-apply paragraph formatting to the newly inserted paragraphs so they
#    match the rest of the document's field list style.
# ---------------------------------------------------------------------------

# The anchor paragraph (originally the blank spacer right after Media) kept
# its own formatting (ind left=1440 hanging=1440) and is now the new blank
# line that sits between "Media:" and "Copy:". Re-number from there.
$p = $mediaPara.Next()              # blank line after Media (unchanged style)
$p = $p.Next()                      # "Copy:  1) ..." paragraph
$p.Range.ParagraphFormat.LeftIndent = 72
$p.Range.ParagraphFormat.FirstLineIndent = -72

$p = $p.Next()                      # blank line
$p.Range.ParagraphFormat.LeftIndent = 72
$p.Range.ParagraphFormat.FirstLineIndent = -72

$p = $p.Next()                      # "2) ..." paragraph
$p.Range.ParagraphFormat.LeftIndent = 72
$p.Range.ParagraphFormat.FirstLineIndent = 0

$p = $p.Next()                      # blank line
$p.Range.ParagraphFormat.LeftIndent = 72
$p.Range.ParagraphFormat.FirstLineIndent = 0

$p = $p.Next()                      # "3) ..." paragraph
$p.Range.ParagraphFormat.LeftIndent = 72
$p.Range.ParagraphFormat.FirstLineIndent = 0

$p = $p.Next()                      # blank line
$p.Range.ParagraphFormat.LeftIndent = 72
$p.Range.ParagraphFormat.FirstLineIndent = 0

$p = $p.Next()                      # "4) ..." paragraph
$p.Range.ParagraphFormat.LeftIndent = 72
$p.Range.ParagraphFormat.FirstLineIndent = 0

$p = $p.Next()                      # blank line
$p.Range.ParagraphFormat.LeftIndent = 72
$p.Range.ParagraphFormat.FirstLineIndent = 0

Write-Host "Edit complete"
